$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 46 (shifts old rows 46:100 down to 47:101)
$ws.Rows.Item(46).Insert()

# Populate the new row 46 with its data
$ws.Cells.Item(46, 1).Value = 7
$ws.Cells.Item(46, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(46, 3).Value = "Ñuble"
$ws.Cells.Item(46, 4).Value = 45128
$ws.Cells.Item(46, 5).Value = 16
$ws.Cells.Item(46, 6).Value = 100112001
$ws.Cells.Item(46, 7).Value = "Berenjena"
$ws.Cells.Item(46, 8).Value = "Sin especificar"
$ws.Cells.Item(46, 9).Value = "Primera"
$ws.Cells.Item(46, 10).Value = 50
$ws.Cells.Item(46, 11).Value = 6500
$ws.Cells.Item(46, 12).Value = 6500
$ws.Cells.Item(46, 13).Value = 6500
$ws.Cells.Item(46, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(46, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(46, 16).Value = 108
$ws.Cells.Item(46, 17).Value = 60
$ws.Cells.Item(46, 18).Value = "Hortaliza"
